$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common values that are identical across every data row (columns B-I)
$origin      = "BA11 5LB"
$destination = "BA11 5AP"
$startAddr   = "81 Knights Maltings, Frome, Frome, BA11 5LB, United Kingdom"
$endAddr     = "55 Tower View, Frome, Frome, BA11 5AP, United Kingdom"
$distText    = 3.0501
$distValue   = 3050.1
$durText     = 8.711666666666668
$durValue    = 522.7

# Per-row Lat/Lng values, keyed by row number (2..30)
$latLng = @{
  2  = @(51.22234, -2.31109)
  3  = @(51.22237, -2.3107)
  4  = @(51.22273, -2.31064)
  5  = @(51.22283, -2.31005)
  6  = @(51.22298, -2.30982)
  7  = @(51.22374, -2.30909)
  8  = @(51.22498, -2.30754)
  9  = @(51.22534, -2.30686)
  10 = @(51.22581, -2.3054)
  11 = @(51.22681, -2.30373)
  12 = @(51.22708, -2.30363)
  13 = @(51.22726, -2.30377)
  14 = @(51.22884, -2.3063)
  15 = @(51.22893, -2.30699)
  16 = @(51.22876, -2.30829)
  17 = @(51.22791, -2.31099)
  18 = @(51.22768, -2.31325)
  19 = @(51.22642, -2.31437)
  20 = @(51.22582, -2.31544)
  21 = @(51.22519, -2.31769)
  22 = @(51.22484, -2.32133)
  23 = @(51.22421, -2.32144)
  24 = @(51.22301, -2.32124)
  25 = @(51.22183, -2.32125)
  26 = @(51.22048, -2.32074)
  27 = @(51.22034, -2.31956)
  28 = @(51.21988, -2.31827)
  29 = @(51.22045, -2.31728)
  30 = @(51.22032, -2.31717)
}

# Update existing rows 2-26 with the new addresses / distances / durations / coordinates
for ($r = 2; $r -le 26; $r++) {
  $ws.Range("B$r").Value = $origin
  $ws.Range("C$r").Value = $destination
  $ws.Range("D$r").Value = $startAddr
  $ws.Range("E$r").Value = $endAddr
  $ws.Range("F$r").Value = $distText
  $ws.Range("G$r").Value = $distValue
  $ws.Range("H$r").Value = $durText
  $ws.Range("I$r").Value = $durValue
  $coords = $latLng[$r]
  $ws.Range("J$r").Value = $coords[0]
  $ws.Range("K$r").Value = $coords[1]
}

# Append new rows 27-30, matching the row-26 formatting (bordered/bold style for column A)
for ($r = 27; $r -le 30; $r++) {
  $ws.Range("A26").Copy($ws.Range("A$r"))
  $ws.Range("A$r").Value = $r - 2
  $ws.Range("B$r").Value = $origin
  $ws.Range("C$r").Value = $destination
  $ws.Range("D$r").Value = $startAddr
  $ws.Range("E$r").Value = $endAddr
  $ws.Range("F$r").Value = $distText
  $ws.Range("G$r").Value = $distValue
  $ws.Range("H$r").Value = $durText
  $ws.Range("I$r").Value = $durValue
  $coords = $latLng[$r]
  $ws.Range("J$r").Value = $coords[0]
  $ws.Range("K$r").Value = $coords[1]
}
